$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = "2024-10-24 00:00:00"
$ws.Range("B83").Value = 73300
$ws.Range("C83").Value = 10269.7
$ws.Range("D83").Value = 9088.23
$ws.Range("E83").Value = 7.1148
